# Automated update 2025-08-27 10:00:10
#
# A new client row ("CORREA IGLESIAS RAMIRO MARCELO") is inserted right
# before "GRANJA VANEGAS MARCELA" (i.e. at row 6) in both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, under the same
# "VACA PANCHI CAROLINA" advisor, with every metric column at 0.
# Everything below shifts down one row, and the trailing "X de N" summary
# row on "VENTAS POR GRUPO" is updated from a denominator of 10 to 11.

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows("6:6").Insert()
$ws1.Range("A6").Value = "VACA PANCHI CAROLINA"
$ws1.Range("B6").Value = "CORREA IGLESIAS RAMIRO MARCELO"
$ws1.Range("C6:R6").Value = 0

# Update the "X de 10" -> "X de 11" summary row, now shifted to row 13.
$ws1.Range("C13:R13").Value = "0 de 11"
$ws1.Range("L13").Value = "1 de 11"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows("6:6").Insert()
$ws2.Range("A6").Value = "VACA PANCHI CAROLINA"
$ws2.Range("B6").Value = "CORREA IGLESIAS RAMIRO MARCELO"
$ws2.Range("C6:G6").Value = 0
